{"js": "const paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nconst lastParagraph = paragraphs.items[paragraphs.items.length - 1];\nlastParagraph.insertParagraph(\"TEST2\", Word.InsertLocation.after);\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# Insert a new paragraph right after the last paragraph in the document body,\n# then fill it with the new text (the paragraph inherits the run/paragraph\n# formatting - en-US language - already in effect at the end of the document).\n$lastPara = $d.Paragraphs.Last\n$lastPara.Range.InsertParagraphAfter()\n\n$newPara = $d.Paragraphs.Last\n$newPara.Range.Text = \"TEST2\"\n"}
